$d = $word.ActiveDocument

$d.Content.Find.Execute("NM2023TMID22194", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NM2023TMID11232", 2)
